$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

$range = $ws.Range("C2:C$lastRow")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45204) {
        $cell.Value2 = 45205
    }
}
